$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$orig = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.169.19'
$ws.Range("D2").Style = $orig
$ws.Range("E2").Value = '  +2.43%  '
$orig = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.305.50'
$ws.Range("D3").Style = $orig
$ws.Range("E3").Value = '  +2.27%  '
$ws.Range("E4").Value = '  +0.00%  '
$orig = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.17'
$ws.Range("D5").Style = $orig
$ws.Range("E5").Value = '  +0.98%  '
$orig = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.65'
$ws.Range("D6").Style = $orig
$ws.Range("E6").Value = '  +4.69%  '
$ws.Range("E7").Value = '  +1.88%  '
$ws.Range("E8").Value = '  -0.03%  '
$orig = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.509'
$ws.Range("D9").Style = $orig
$ws.Range("E9").Value = '  +3.48%  '
$orig = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.07'
$ws.Range("D10").Style = $orig
$ws.Range("E10").Value = '  +4.08%  '
$orig = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0797'
$ws.Range("D11").Style = $orig
$ws.Range("E11").Value = '  +1.72%  '
$orig = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '49.18'
$ws.Range("D12").Style = $orig
$ws.Range("E12").Value = '  +1.78%  '
$ws.Range("E13").Value = '  +3.45%  '
$orig = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.88'
$ws.Range("D14").Style = $orig
$ws.Range("E14").Value = '  +16.87%  '
$orig = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.80'
$ws.Range("D15").Style = $orig
$ws.Range("E15").Value = '  +2.71%  '
$orig = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.667.09'
$ws.Range("D16").Style = $orig
$ws.Range("E16").Value = '  +2.56%  '
$orig = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.311.22'
$ws.Range("D17").Style = $orig
$ws.Range("E17").Value = '  +1.07%  '
$ws.Range("E18").Value = '  +4.86%  '
$orig = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.039.86'
$ws.Range("D19").Style = $orig
$ws.Range("E19").Value = '  +2.29%  '
$orig = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.06'
$ws.Range("D20").Style = $orig
$ws.Range("E20").Value = '  +6.09%  '
$orig = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0908'
$ws.Range("D21").Style = $orig
$ws.Range("E21").Value = '  +2.28%  '
$orig = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.10'
$ws.Range("D22").Style = $orig
$ws.Range("E22").Value = '  +2.22%  '
$orig = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.86'
$ws.Range("D23").Style = $orig
$ws.Range("E23").Value = '  +2.48%  '
$orig = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '237.21'
$ws.Range("D24").Style = $orig
$ws.Range("E24").Value = '  +2.04%  '
$orig = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.08'
$ws.Range("D25").Style = $orig
$ws.Range("E25").Value = '  +8.59%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  +0.18%  '
$orig = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.66'
$ws.Range("D28").Style = $orig
$ws.Range("E28").Value = '  +3.62%  '
$orig = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").Style = $orig
$ws.Range("E29").Value = '  +7.21%  '
$orig = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.02'
$ws.Range("D30").Style = $orig
$ws.Range("E30").Value = '  +0.02%  '
$orig = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.95'
$ws.Range("D31").Style = $orig
$ws.Range("E31").Value = '  +1.02%  '
$orig = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.19'
$ws.Range("D32").Style = $orig
$ws.Range("E32").Value = '  +1.73%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("E34").Value = '  +1.85%  '
$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$orig = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.56'
$ws.Range("D35").Style = $orig
$ws.Range("E35").Value = '  +5.25%  '
$ws.Range("B36").Value = 'WEMIXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$orig = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.43'
$ws.Range("D36").Style = $orig
$ws.Range("E36").Value = '  +4.62%  '
$orig = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.09'
$ws.Range("D37").Style = $orig
$ws.Range("E37").Value = '  +7.47%  '
$orig = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0700'
$ws.Range("D38").Style = $orig
$ws.Range("E38").Value = '  +1.37%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$orig = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.102'
$ws.Range("D39").Style = $orig
$ws.Range("E39").Value = '  +4.28%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$orig = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.84'
$ws.Range("D40").Style = $orig
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$orig = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.80'
$ws.Range("D41").Style = $orig
$ws.Range("E41").Value = '  +5.22%  '
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("E43").Value = '  -2.64%  '
$orig = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.999.41'
$ws.Range("D44").Style = $orig
$ws.Range("E44").Value = '  +3.17%  '
$ws.Range("E45").Value = '  +3.17%  '
$orig = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.02'
$ws.Range("D46").Style = $orig
$ws.Range("E46").Value = '  +5.71%  '
$orig = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.85'
$ws.Range("D47").Style = $orig
$ws.Range("E47").Value = '  +3.05%  '
$ws.Range("E48").Value = '  +4.14%  '
$orig = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.33'
$ws.Range("D49").Style = $orig
$ws.Range("E49").Value = '  +4.33%  '
$orig = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.532.18'
$ws.Range("D50").Style = $orig
$ws.Range("E50").Value = '  +2.08%  '
$ws.Range("E51").Value = '  +3.04%  '
